$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 855
$ws.Range("C3").Value = 855
$ws.Range("C4").Value = 855
$ws.Range("C5").Value = 855
$ws.Range("C6").Value = 855
$ws.Range("C7").Value = 491
$ws.Range("C8").Value = 465
$ws.Range("C9").Value = 531
$ws.Range("C10").Value = 504
$ws.Range("C11").Value = 494
$ws.Range("C12").Value = 679
$ws.Range("C13").Value = 471
$ws.Range("C14").Value = 376
$ws.Range("C15").Value = 413
